# Updates crypto price/volume data per commit:
# "Updated cryptos list on Wed Feb  7 20:58:27 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.129.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.437.97"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.96"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.34"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.67"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.812.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.435.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.834"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.154.79"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.30"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.55"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.89"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.21"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.79"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +15.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.60"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.89%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.92"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "131.69"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +22.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.49"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.95"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.29"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.954.51"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.36%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.668.40"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.63"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +7.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.36"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.30"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.52%  "
